$wb = $excel.ActiveWorkbook

# This script applies the numeric refresh captured in the scheduled-runner commit.
# For each affected worksheet we push the updated currentAveragePrice/LevePrice/LeveProfit
# figures (columns H-N) into the specific rows that changed, then clear any cells whose
# value went from a concrete number to "no cached value" in the new snapshot.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7579.7144
$ws.Range("I74").Value = 5667.3335
$ws.Range("K74").Value = 5667.3335
$ws.Range("M74").Value = -4731.3335
$ws.Range("H77").Value = 7579.7144
$ws.Range("I77").Value = 5667.3335
$ws.Range("K77").Value = 28336.6675
$ws.Range("M77").Value = -23656.6675
$ws.Range("H132").Value = 24393546
$ws.Range("I132").Value = 26319256
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 78957768
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -78955238
$ws.Range("N132").Value = -8660
$ws.Range("H135").Value = 919.2381
$ws.Range("J135").Value = 3299.4
$ws.Range("L135").Value = 29694.6
$ws.Range("N135").Value = -34764.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4065.011
$ws.Range("I32").Value = 2464.5066
$ws.Range("J32").Value = 12280.934
$ws.Range("K32").Value = 2464.5066
$ws.Range("L32").Value = 12280.934
$ws.Range("M32").Value = -2177.5066
$ws.Range("N32").Value = -12854.934
$ws.Range("H92").Value = 84699
$ws.Range("J92").Value = 84699
$ws.Range("L92").Value = 84699
$ws.Range("N92").Value = -89691
$ws.Range("H122").Value = 1044557.94
$ws.Range("I122").Value = 2703.1428
$ws.Range("J122").Value = 3475552.5
$ws.Range("K122").Value = 8109.428400000001
$ws.Range("L122").Value = 10426657.5
$ws.Range("M122").Value = -5659.428400000001
$ws.Range("N122").Value = -10431557.5
$ws.Range("H132").Value = 1555.4584
$ws.Range("I132").Value = 1272.2354
$ws.Range("K132").Value = 3816.7062
$ws.Range("M132").Value = -1286.7062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3710847
$ws.Range("I94").Value = 5263479
$ws.Range("K94").Value = 5263479
$ws.Range("M94").Value = -5263028
$ws.Range("H99").Value = 7574250
$ws.Range("I99").Value = 11068157
$ws.Range("J99").Value = 4117.5
$ws.Range("K99").Value = 11068157
$ws.Range("L99").Value = 4117.5
$ws.Range("M99").Value = -11066659
$ws.Range("N99").Value = -7113.5
$ws.Range("H105").Value = 3909707.2
$ws.Range("I105").Value = 4170154.2
$ws.Range("K105").Value = 4170154.2
$ws.Range("M105").Value = -4168407.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1719.1154
$ws.Range("I58").Value = 1391.75
$ws.Range("J58").Value = 2810.3333
$ws.Range("K58").Value = 1391.75
$ws.Range("L58").Value = 2810.3333
$ws.Range("M58").Value = -1188.75
$ws.Range("N58").Value = -3216.3333
$ws.Range("H69").Value = 41244.184
$ws.Range("I69").Value = 18556.715
$ws.Range("K69").Value = 18556.715
$ws.Range("M69").Value = -17807.715
$ws.Range("H72").Value = 41244.184
$ws.Range("I72").Value = 18556.715
$ws.Range("K72").Value = 55670.145
$ws.Range("M72").Value = -51926.145
$ws.Range("H105").Value = 2958.7144
$ws.Range("I105").Value = 2946.0908
$ws.Range("J105").Value = 3005
$ws.Range("K105").Value = 2946.0908
$ws.Range("L105").Value = 3005
$ws.Range("M105").Value = -1199.0908
$ws.Range("N105").Value = -6499
$ws.Range("H132").Value = 15849.203
$ws.Range("I132").Value = 1746.3914
$ws.Range("J132").Value = 51889.723
$ws.Range("K132").Value = 5239.174199999999
$ws.Range("L132").Value = 155669.169
$ws.Range("M132").Value = -2709.174199999999
$ws.Range("N132").Value = -160729.169
$ws.Range("H133").Value = 40217.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 40217.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 40217.332
$ws.Range("N133").Value = -45277.332
$ws.Range("H136").Value = 1719.1154
$ws.Range("I136").Value = 1391.75
$ws.Range("J136").Value = 2810.3333
$ws.Range("K136").Value = 4175.25
$ws.Range("L136").Value = 8430.999899999999
$ws.Range("M136").Value = -1625.25
$ws.Range("N136").Value = -13530.9999
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 499.25
$ws.Range("J75").Value = 594.6667
$ws.Range("L75").Value = 1784.0001
$ws.Range("N75").Value = -3780.0001
$ws.Range("H78").Value = 499.25
$ws.Range("J78").Value = 594.6667
$ws.Range("L78").Value = 5352.0003
$ws.Range("N78").Value = -15336.0003
$ws.Range("H80").Value = 4621
$ws.Range("J80").Value = 2828
$ws.Range("L80").Value = 8484
$ws.Range("N80").Value = -10356
$ws.Range("H83").Value = 4621
$ws.Range("J83").Value = 2828
$ws.Range("L83").Value = 25452
$ws.Range("N83").Value = -34812
$ws.Range("H98").Value = 719.25
$ws.Range("I98").Value = 620.6667
$ws.Range("J98").Value = 817.8333
$ws.Range("K98").Value = 1862.0001
$ws.Range("L98").Value = 2453.4999
$ws.Range("M98").Value = -364.0001
$ws.Range("N98").Value = -5449.4999
$ws.Range("H103").Value = 150.625
$ws.Range("J103").Value = 163.33333
$ws.Range("L103").Value = 489.99999
$ws.Range("N103").Value = -2247.99999
$ws.Range("H114").Value = 1187.2142
$ws.Range("J114").Value = 1140.0769
$ws.Range("L114").Value = 3420.2307
$ws.Range("N114").Value = -9928.2307
$ws.Range("H121").Value = 662.75
$ws.Range("I121").Value = 342.25
$ws.Range("J121").Value = 983.25
$ws.Range("K121").Value = 1026.75
$ws.Range("L121").Value = 2949.75
$ws.Range("M121").Value = 283.25
$ws.Range("N121").Value = -5569.75
$ws.Range("H129").Value = 992.05884
$ws.Range("I129").Value = 914.3333
$ws.Range("J129").Value = 1178.6
$ws.Range("K129").Value = 2742.9999
$ws.Range("L129").Value = 3535.8
$ws.Range("M129").Value = 2257.0001
$ws.Range("N129").Value = -13535.8
$ws.Range("H131").Value = 7864607.5
$ws.Range("I131").Value = 4387492.5
$ws.Range("K131").Value = 13162477.5
$ws.Range("M131").Value = -13157437.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 6900
$ws.Range("I20").Value = 6800
$ws.Range("J20").Value = 7000
$ws.Range("K20").Value = 6800
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = -6555
$ws.Range("N20").Value = -7490
$ws.Range("H92").Value = 17262.5
$ws.Range("J92").Value = 17262.5
$ws.Range("L92").Value = 17262.5
$ws.Range("N92").Value = -21006.5
$ws.Range("H102").Value = 6748690.5
$ws.Range("I102").Value = 7938449.5
$ws.Range("J102").Value = 3972586.5
$ws.Range("K102").Value = 7938449.5
$ws.Range("L102").Value = 3972586.5
$ws.Range("M102").Value = -7936827.5
$ws.Range("N102").Value = -3975830.5
$ws.Range("H116").Value = 100214
$ws.Range("J116").Value = 100214
$ws.Range("L116").Value = 100214
$ws.Range("N116").Value = -109392
$ws.Range("H117").Value = 29999.572
$ws.Range("J117").Value = 29999.572
$ws.Range("L117").Value = 29999.572
$ws.Range("N117").Value = -36883.572
$ws.Range("H119").Value = 90000
$ws.Range("J119").Value = 90000
$ws.Range("L119").Value = 90000
$ws.Range("N119").Value = -99676
$ws.Range("H120").Value = 38750
$ws.Range("J120").Value = 38750
$ws.Range("L120").Value = 38750
$ws.Range("N120").Value = -48426
$ws.Range("H122").Value = 735364.4
$ws.Range("I122").Value = 836651.4
$ws.Range("K122").Value = 2509954.2
$ws.Range("M122").Value = -2507504.2
$ws.Range("H126").Value = 8589182
$ws.Range("I126").Value = 4548564
$ws.Range("J126").Value = 16670416
$ws.Range("K126").Value = 13645692
$ws.Range("L126").Value = 50011248
$ws.Range("M126").Value = -13643222
$ws.Range("N126").Value = -50016188

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 18000
$ws.Range("I5").Value = 18000
$ws.Range("K5").Value = 18000
$ws.Range("M5").Value = -17887
$ws.Range("H29").Value = 28000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("H40").Value = 14971
$ws.Range("I40").Value = 14943
$ws.Range("K40").Value = 14943
$ws.Range("M40").Value = -14807
$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25386
$ws.Range("H68").Value = 4596.4
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 4596.4
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H136").Value = 30431.63
$ws.Range("I136").Value = 55142
$ws.Range("K136").Value = 165426
$ws.Range("M136").Value = -162876
$ws.Range("M29").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 18748.5
$ws.Range("I37").Value = 18748.5
$ws.Range("K37").Value = 18748.5
$ws.Range("M37").Value = -18545.5
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("H107").Value = 83334170
$ws.Range("I107").Value = 142857800
$ws.Range("J107").Value = 1099
$ws.Range("K107").Value = 428573400
$ws.Range("L107").Value = 3297
$ws.Range("M107").Value = -428571480
$ws.Range("N107").Value = -7137
$ws.Range("H113").Value = 870.8929000000001
$ws.Range("I113").Value = 354.88235
$ws.Range("K113").Value = 1064.64705
$ws.Range("M113").Value = 1105.35295
$ws.Range("H132").Value = 32291644
$ws.Range("I132").Value = 41672380
$ws.Range("J132").Value = 129116.43
$ws.Range("K132").Value = 125017140
$ws.Range("L132").Value = 387349.29
$ws.Range("M132").Value = -125014610
$ws.Range("N132").Value = -392409.29
$ws.Range("M42").ClearContents()
